# Adding the changes we made on may 9th
#
# Inserts 11 new data rows at the top of the table (rows 2-12), pushing the
# existing data down by 11 rows, then removes what is now the trailing
# (previously-last) row so the sheet ends up with rows 2-31 of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data (rows 2:21) down by 11 rows to make room for the
# new rows, i.e. insert 11 blank rows starting at row 2.
$ws.Rows("2:12").Insert()

# Insert() copies the formatting of the row above by default; the new rows
# should be plain (unstyled) like the rest of the numeric data, so strip
# whatever formatting just got carried down.
$ws.Range("A2:C12").ClearFormats()

# New data for rows 2-12.
$newData = @(
    @(0.308792382478714, -0.2591595947742462, 0.3060434758663177),
    @(0.197004035115242, 0.0836885422468185, 0.0430659987032413),
    @(0.0916297882795333, 0.034972034394741, 0.0708603709936142),
    @(0.0360410511493682, 0.06902777403593061, -0.0134390350431203),
    @(-0.0181732401251792, 0.0655152946710586, -0.0574213340878486),
    @(-0.1020144969224929, -0.0171042270958423, -0.0577267669141292),
    @(-0.0662788823246955, -0.0591012127697467, -0.0519235469400882),
    @(-0.0091629782691597, 0.016951510682702, 0.0032070425804704),
    @(-0.0189368221908807, 0.0145080499351024, 0.0222965814173221),
    @(-0.0164933614432811, -0.0006108652451075, -0.0030543261673301),
    @(-0.0395535230636596, -0.0021380283869802, 0.0117591563612222)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = 2 + $i
    $vals = $newData[$i]
    $ws.Range("A$row").Value = $vals[0]
    $ws.Range("B$row").Value = $vals[1]
    $ws.Range("C$row").Value = $vals[2]
}

# The insert shifted the old trailing row (formerly row 21) down to row 32;
# that row is dropped entirely in the final layout, so delete it.
$ws.Rows("32:32").Delete()
